$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item(1)

# New rows of data appended to the bottom of the list (A137:B151).
# Cells are written in the same order the source data was entered so the
# shared-string table ends up populated in the same sequence.
$ws.Cells.Item(137, 1).Value = "giặt sấy"
$ws.Cells.Item(140, 1).Value = "nhà thuốc"
$ws.Cells.Item(138, 1).Value = "nhà thuốc tư nhân"
$ws.Cells.Item(141, 1).Value = "vựa khô đường đậu gia vị "
$ws.Cells.Item(142, 1).Value = "tã sữa tổng hợp"
$ws.Cells.Item(139, 1).Value = "nhà thuốc tây"
$ws.Cells.Item(143, 1).Value = "lan ,mỹ phẩm"
$ws.Cells.Item(143, 2).Value = "lan"
$ws.Cells.Item(144, 1).Value = "chuyên bán sỉ lẻ nước ngọt"
$ws.Cells.Item(145, 1).Value = "cưa hàng gia dụng"
$ws.Cells.Item(146, 1).Value = "phụ liệu tóc nail"
$ws.Cells.Item(147, 1).Value = "nhận sửa quần áo"
$ws.Cells.Item(148, 1).Value = "chuyên cung cấp các loại sữa"
$ws.Cells.Item(149, 1).Value = "tạp hoa gia dụng"
$ws.Cells.Item(150, 1).Value = "xe gắn máy"
$ws.Cells.Item(151, 1).Value = "ehome"

# Update the selection to reflect the new last cell in the list
$ws.Activate()
$ws.Range("A151").Select()
